# 1st commit - Data-driven - Adding parameterization
#
# Inserts a new "test_suite" worksheet (TCID / Runmode table) as the first
# sheet in the workbook, driving AddCustomerTest / OpenAccountTest from it.

$wb = $excel.ActiveWorkbook

# Insert the new "test_suite" sheet before the current first sheet so it
# becomes sheet1 (AddCustomerTest/OpenAccountTest shift down one slot).
$testSuite = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$testSuite.Name = "test_suite"

# Header row.
$testSuite.Range("A1").Value = "TCID"
$testSuite.Range("B1").Value = "Runmode"

# TCID column first (keeps shared-string insertion order: TCID, Runmode,
# BankManagerLoginTest, AddCustomerTest, OpenAccountTest, then Y).
$testSuite.Range("A2").Value = "BankManagerLoginTest"
$testSuite.Range("A3").Value = "AddCustomerTest"
$testSuite.Range("A4").Value = "OpenAccountTest"

# Runmode column.
$testSuite.Range("B2").Value = "Y"
$testSuite.Range("B3").Value = "Y"
$testSuite.Range("B4").Value = "Y"

# Best-fit column A to the longest label ("BankManagerLoginTest").
$testSuite.Columns.Item(1).ColumnWidth = 21

# The AddCustomerTest sheet is no longer the active tab; its cursor moves
# to N1 and it keeps whatever data it already had.
$addCustomer = $wb.Worksheets.Item("AddCustomerTest")
$addCustomer.Range("N1").Select()

# test_suite becomes (and stays) the active tab, cursor parked at A8.
$testSuite.Range("A8").Select()
